$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price / link / volume(1h) data (GitHub Actions scrape update)
$ws.Range('D2').Value = '29.068.65'
$ws.Range('E2').Value = '  +2.77%  '
$ws.Range('D3').Value = '1.582.36'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').Value = '  -0.57%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  +6.98%  '
$ws.Range('E7').Value = '  -0.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '25.58'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.03%  '
$ws.Range('E9').Value = '  +2.79%  '
$ws.Range('E10').Value = '  +1.88%  '
$ws.Range('E11').Value = '  +1.24%  '
$ws.Range('D12').Value = '1.808.65'
$ws.Range('E12').Value = '  +1.96%  '
$ws.Range('D13').Value = '1.612.76'
$ws.Range('E13').Value = '  +3.75%  '
$ws.Range('D14').Value = '29.113.73'
$ws.Range('E14').Value = '  +2.92%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.522'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.67%  '
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.70'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.50'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '238.12'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.03%  '
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range('E20').Value = '  +3.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.996'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.19'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.17%  '
$ws.Range('E24').Value = '  +4.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.03'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.66%  '
$ws.Range('E26').Value = '  +5.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.13'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.50%  '
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('E30').Value = '  -0.61%  '
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('E32').Value = '  +1.38%  '
$ws.Range('D33').Value = '1.421.03'
$ws.Range('E33').Value = '  +2.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.05'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.32%  '
$ws.Range('E35').Value = '  -1.64%  '
$ws.Range('E36').Value = '  +1.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.79'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.20%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.29'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.99%  '
$ws.Range('E39').Value = '  +1.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.525'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.09%  '
$ws.Range('E41').Value = '  +0.64%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '53.30'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +25.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.996'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.60%  '
$ws.Range('E44').Value = '  +1.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0461'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.61'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.65%  '
$ws.Range('E47').Value = '  -1.58%  '
$ws.Range('D48').Value = '1.720.18'
$ws.Range('E48').Value = '  +1.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.850'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.24%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '85.59'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0513'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.83%  '
